$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Origin")

# Remove rows 4 through 11 content (only rows 1-3 remain: header + 2 data rows)
$ws.Range("A4:E11").ClearContents() | Out-Null

# Update row 2 / row 3 data (order matches original authoring sequence)
$ws.Range("C2").Value = "AutomatedTest01"
$ws.Range("D3").Value = "Sinopharm Group Industries"
$ws.Range("C3").Value = "AutomatedTest02"
$ws.Range("D2").Value = "Takeda Pharmaceutical Industries"
$ws.Range("E2").Value = "Japan"
$ws.Range("E3").Value = "Tiongkok"

# Update selection to D2
$ws.Range("D2").Select() | Out-Null
